# Applies the "Updated cryptos list" data refresh (Wed Aug 14 18:35:31 UTC 2024).
# For each changed cell we temporarily force a Text number format so that
# numeric-looking strings (e.g. "1.00", "0.0230") are not silently coerced
# into numbers (which would drop the formatting), then restore the default
# "Normal" style so the cell format matches the original workbook exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '59.165.72'
Set-TextValue 'E2' '  -2.80%  '
Set-TextValue 'D3' '2.651.75'
Set-TextValue 'E3' '  -1.18%  '
Set-TextValue 'E4' '  +0.12%  '
Set-TextValue 'D5' '524.14'
Set-TextValue 'E5' '  +0.55%  '
Set-TextValue 'D6' '144.84'
Set-TextValue 'E6' '  -1.45%  '
Set-TextValue 'D7' '1.00'
Set-TextValue 'E7' '  +0.19%  '
Set-TextValue 'D8' '0.572'
Set-TextValue 'E8' '  -1.18%  '
Set-TextValue 'D9' '7.04'
Set-TextValue 'E9' '  +9.39%  '
Set-TextValue 'E10' '  -3.18%  '
Set-TextValue 'E11' '  -2.02%  '
Set-TextValue 'E12' '  +1.71%  '
Set-TextValue 'D13' '3.122.63'
Set-TextValue 'E13' '  -1.13%  '
Set-TextValue 'D14' '59.173.39'
Set-TextValue 'E14' '  -2.71%  '
Set-TextValue 'D15' '21.10'
Set-TextValue 'E15' '  -1.41%  '
Set-TextValue 'B16' 'WrappedEther'
Set-TextValue 'C16' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D16' '2.670.32'
Set-TextValue 'E16' '  -3.10%  '
Set-TextValue 'B17' 'ShibaInu'
Set-TextValue 'C17' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue 'D17' '0.0000136'
Set-TextValue 'E17' '  -1.87%  '
Set-TextValue 'D18' '340.56'
Set-TextValue 'E18' '  -3.49%  '
Set-TextValue 'D19' '4.38'
Set-TextValue 'E19' '  -4.13%  '
Set-TextValue 'D20' '10.38'
Set-TextValue 'E20' '  -1.72%  '
Set-TextValue 'D21' '6.38'
Set-TextValue 'E21' '  +0.31%  '
Set-TextValue 'D22' '1.00'
Set-TextValue 'E22' '  +0.32%  '
Set-TextValue 'D23' '64.41'
Set-TextValue 'E23' '  +2.25%  '
Set-TextValue 'E24' '  -0.87%  '
Set-TextValue 'E25' '  -1.49%  '
Set-TextValue 'E26' '  +0.40%  '
Set-TextValue 'E27' '  -1.92%  '
Set-TextValue 'D28' '7.13'
Set-TextValue 'E28' '  -2.49%  '
Set-TextValue 'E29' '  -2.39%  '
Set-TextValue 'D30' '0.999'
Set-TextValue 'E30' '  +0.02%  '
Set-TextValue 'E31' '  +0.03%  '
Set-TextValue 'D32' '18.87'
Set-TextValue 'E32' '  -1.28%  '
Set-TextValue 'D33' '149.29'
Set-TextValue 'E33' '  -0.16%  '
Set-TextValue 'E34' '  -3.41%  '
Set-TextValue 'E35' '  -3.15%  '
Set-TextValue 'D36' '0.895'
Set-TextValue 'E36' '  -5.71%  '
Set-TextValue 'D37' '0.878'
Set-TextValue 'E37' '  +0.25%  '
Set-TextValue 'D38' '36.75'
Set-TextValue 'E38' '  +0.16%  '
Set-TextValue 'E39' '  -5.73%  '
Set-TextValue 'E40' '  -3.22%  '
Set-TextValue 'D41' '0.618'
Set-TextValue 'E41' '  +0.94%  '
Set-TextValue 'D42' '20.04'
Set-TextValue 'E42' '  -0.20%  '
Set-TextValue 'E43' '  +0.15%  '
Set-TextValue 'D44' '275.01'
Set-TextValue 'E44' '  -3.29%  '
Set-TextValue 'E45' '  -2.09%  '
Set-TextValue 'B46' 'WhiteBITCoin'
Set-TextValue 'C46' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D46' '10.66'
Set-TextValue 'E46' '  +1.87%  '
Set-TextValue 'B47' 'Hedera'
Set-TextValue 'C47' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue 'D47' '0.0534'
Set-TextValue 'E47' '  -1.54%  '
Set-TextValue 'B48' 'RenderToken'
Set-TextValue 'C48' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D48' '4.81'
Set-TextValue 'E48' '  -1.34%  '
Set-TextValue 'B49' 'Maker'
Set-TextValue 'C49' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue 'D49' '2.033.99'
Set-TextValue 'E49' '  -4.68%  '
Set-TextValue 'D50' '0.0230'
Set-TextValue 'E50' '  -2.73%  '
Set-TextValue 'D51' '18.94'
Set-TextValue 'E51' '  -1.00%  '

Write-Output "Applied 96 cell updates"
